$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 55
$ws1.Range("F3").Value = 102
$ws1.Range("F4").Value = 7345
$ws1.Range("F7").Value = 3862
$ws1.Range("F8").Value = 312
$ws1.Range("F9").Value = 551
$ws1.Range("F12").Value = 111

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 55
$ws4.Range("F3").Value = 102
$ws4.Range("F5").Value = 7345
$ws4.Range("F9").Value = 3862
$ws4.Range("F10").Value = 312
$ws4.Range("F11").Value = 551
$ws4.Range("F14").Value = 111
